# Update gh-pages to output generated at 456a3b4
# Applies updated "want to go" counts / prices / address info to the
# 展览 (Exhibition) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # F column: 想去人数 (number of people interested)
    $ws.Range("F4").Value = 52
    $ws.Range("F7").Value = 151

    # D9: 地点 (location) changed
    $ws.Range("D9").Value = "真君路888号 南昌华侨城玩美公园"
    $ws.Range("F9").Value = 24

    $ws.Range("F10").Value = 250
    $ws.Range("F15").Value = 658
    $ws.Range("F17").Value = 484
    $ws.Range("F18").Value = 415

    $ws.Range("F23").Value = 1231
    $ws.Range("G23").Value = 65

    $ws.Range("F24").Value = 2894
    $ws.Range("F25").Value = 23

    $ws.Range("F28").Value = 60
    $ws.Range("F29").Value = 1623
    $ws.Range("F32").Value = 14
    $ws.Range("F36").Value = 609
    $ws.Range("F37").Value = 423
    $ws.Range("F38").Value = 5
}

# F27 differs slightly between the two sheets
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F27").Value = 550

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F27").Value = 551
